$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of metric data (row 9)
$ws.Cells.Item(9, 1).Value = "2025-04-28 10:55:49"
$ws.Cells.Item(9, 2).Value = 223
